$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells remain text (avoid Excel auto-numeric coercion),
# then restore default (unstyled) formatting so the saved XML matches the
# original un-styled inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.967.33'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '1.713.45'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '318.01'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '0.3971'
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("D8").Value = '0.4121'
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("D9").Value = '1.531'
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("D10").Value = '1.002'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = '53.66'
$ws.Range("E11").Value = '  +3.65%  '
$ws.Range("D12").Value = '0.08957'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").Value = '7.752'
$ws.Range("E13").Value = '  +7.28%  '
$ws.Range("D14").Value = '24.81'
$ws.Range("E14").Value = '  +5.87%  '
$ws.Range("D15").Value = '8.177'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '0.00001374'
$ws.Range("E16").Value = '  +4.28%  '
$ws.Range("D17").Value = '1.691.93'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = '100.52'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").Value = '0.07147'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").Value = '20.18'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").Value = '7.525'
$ws.Range("E21").Value = '  +6.20%  '
$ws.Range("D22").Value = '1.006'
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").Value = '14.57'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("D24").Value = '24.965.93'
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").Value = '3.145'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").Value = '2.331'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = '23.19'
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("D28").Value = '9.325'
$ws.Range("E28").Value = '  +24.44%  '
$ws.Range("D29").Value = '165.54'
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").Value = '140.59'
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("D31").Value = '5.229'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '7.925'
$ws.Range("E32").Value = '  +11.18%  '
$ws.Range("D33").Value = '0.09058'
$ws.Range("E33").Value = '  +4.72%  '
$ws.Range("D34").Value = '1.878.32'
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").Value = '1.092'
$ws.Range("E35").Value = '  +1.09%  '
$ws.Range("D36").Value = '0.03020'
$ws.Range("E36").Value = '  +10.80%  '
$ws.Range("D37").Value = '0.2815'
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("D38").Value = '11.16'
$ws.Range("E38").Value = '  -3.65%  '
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("D40").Value = '14.62'
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").Value = '0.09315'
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("D42").Value = '0.8126'
$ws.Range("E42").Value = '  +6.03%  '
$ws.Range("D43").Value = '1.489'
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("D44").Value = '16.74'
$ws.Range("E44").Value = '  +7.04%  '
$ws.Range("D45").Value = '0.7406'
$ws.Range("E45").Value = '  +3.10%  '
$ws.Range("D46").Value = '2.652'
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = '4.273'
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").Value = '1.353'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").Value = '140.96'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '93.81'
$ws.Range("E51").Value = '  +4.52%  '

$ws.Range("D2:D51").Style = "Normal"

